$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B6: change from text to a real number
$ws.Range("B6").Value = 67890543

# Add a new row 7 with Elton John's data
$ws.Range("A7").Value = "Elton John"

# B7 keeps the phone number as text (matches the source data)
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "3456789"

$ws.Range("C7").Value = "eltonjohn@gmail.com"
